$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every data cell in this sheet is stored as text (inline string), even when
# the text looks numeric (e.g. "293.06", "0.0790"). Excel's Range.Value setter
# auto-coerces numeric-looking strings to real numbers, which would both change
# the cell type and silently drop formatting like trailing zeros. Forcing the
# cell to Text format ("@") right before the write keeps the literal string,
# and ClearFormats() right after drops the now-unneeded explicit number format
# so the cell's style index is left exactly as it was (unstyled / style 0).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "40.082.51"
Set-TextValue $ws.Range("E2") "  +1.83%  "
Set-TextValue $ws.Range("D3") "2.236.70"
Set-TextValue $ws.Range("E3") "  +0.77%  "
Set-TextValue $ws.Range("E4") "  +0.17%  "
Set-TextValue $ws.Range("D5") "293.06"
Set-TextValue $ws.Range("E5") "  -1.14%  "
Set-TextValue $ws.Range("D6") "87.11"
Set-TextValue $ws.Range("E6") "  +4.72%  "
Set-TextValue $ws.Range("D7") "0.516"
Set-TextValue $ws.Range("E7") "  +1.10%  "
Set-TextValue $ws.Range("E8") "  +0.07%  "
Set-TextValue $ws.Range("E9") "  +1.35%  "
Set-TextValue $ws.Range("D10") "31.27"
Set-TextValue $ws.Range("E10") "  +7.35%  "
Set-TextValue $ws.Range("D11") "0.0790"
Set-TextValue $ws.Range("E11") "  +2.02%  "
Set-TextValue $ws.Range("D12") "46.97"
Set-TextValue $ws.Range("E12") "  -1.92%  "
Set-TextValue $ws.Range("E13") "  +1.34%  "
Set-TextValue $ws.Range("D14") "6.42"
Set-TextValue $ws.Range("E14") "  +1.82%  "
Set-TextValue $ws.Range("D15") "2.586.37"
Set-TextValue $ws.Range("E15") "  +0.86%  "
Set-TextValue $ws.Range("D16") "14.10"
Set-TextValue $ws.Range("E16") "  -0.06%  "
Set-TextValue $ws.Range("D17") "2.239.68"
Set-TextValue $ws.Range("E17") "  +1.71%  "
Set-TextValue $ws.Range("E18") "  +2.62%  "
Set-TextValue $ws.Range("D19") "40.015.57"
Set-TextValue $ws.Range("E19") "  +1.93%  "
Set-TextValue $ws.Range("E20") "  +1.81%  "
Set-TextValue $ws.Range("D21") "11.29"
Set-TextValue $ws.Range("E21") "  +9.86%  "
Set-TextValue $ws.Range("D22") "5.83"
Set-TextValue $ws.Range("E22") "  +2.13%  "
Set-TextValue $ws.Range("D23") "65.70"
Set-TextValue $ws.Range("E23") "  +1.03%  "
Set-TextValue $ws.Range("D24") "236.05"
Set-TextValue $ws.Range("E24") "  +3.57%  "
Set-TextValue $ws.Range("D25") "0.999"
Set-TextValue $ws.Range("E25") "  -0.17%  "
Set-TextValue $ws.Range("D27") "1.85"
Set-TextValue $ws.Range("E27") "  +2.41%  "
Set-TextValue $ws.Range("D28") "22.95"
Set-TextValue $ws.Range("E28") "  +1.73%  "
Set-TextValue $ws.Range("E29") "  +2.86%  "
Set-TextValue $ws.Range("D30") "9.33"
Set-TextValue $ws.Range("E30") "  +2.45%  "
Set-TextValue $ws.Range("D31") "33.39"
Set-TextValue $ws.Range("E31") "  +4.52%  "
Set-TextValue $ws.Range("D32") "151.51"
Set-TextValue $ws.Range("E32") "  +1.14%  "
Set-TextValue $ws.Range("E33") "  -0.03%  "
Set-TextValue $ws.Range("D34") "4.93"
Set-TextValue $ws.Range("E34") "  +1.72%  "
Set-TextValue $ws.Range("D35") "0.0721"
Set-TextValue $ws.Range("E35") "  +3.96%  "
Set-TextValue $ws.Range("E36") "  +2.27%  "
Set-TextValue $ws.Range("D37") "16.33"
Set-TextValue $ws.Range("E37") "  +7.63%  "
Set-TextValue $ws.Range("D38") "2.82"
Set-TextValue $ws.Range("E38") "  +6.89%  "
Set-TextValue $ws.Range("E39") "  +2.20%  "
Set-TextValue $ws.Range("D40") "0.1000"
Set-TextValue $ws.Range("E40") "  +3.87%  "
Set-TextValue $ws.Range("D41") "1.72"
Set-TextValue $ws.Range("E41") "  +4.74%  "
Set-TextValue $ws.Range("D42") "3.84"
Set-TextValue $ws.Range("E42") "  +5.44%  "
Set-TextValue $ws.Range("D43") "2.067.50"
Set-TextValue $ws.Range("E43") "  +8.32%  "
Set-TextValue $ws.Range("D44") "18.16"
Set-TextValue $ws.Range("E44") "  +13.39%  "
Set-TextValue $ws.Range("D45") "0.0270"
Set-TextValue $ws.Range("E45") "  +4.19%  "
Set-TextValue $ws.Range("E46") "  +4.31%  "
Set-TextValue $ws.Range("D47") "9.81"
Set-TextValue $ws.Range("E47") "  +9.79%  "
Set-TextValue $ws.Range("D48") "2.60"
Set-TextValue $ws.Range("E48") "  -0.72%  "
Set-TextValue $ws.Range("B49") "RocketPoolETH"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue $ws.Range("D49") "2.445.65"
Set-TextValue $ws.Range("E49") "  +0.49%  "
Set-TextValue $ws.Range("B50") "BitcoinSV"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue $ws.Range("D50") "72.08"
Set-TextValue $ws.Range("E50") "  +2.38%  "
Set-TextValue $ws.Range("D51") "89.35"
Set-TextValue $ws.Range("E51") "  +2.79%  "
